$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.15828275680542
$ws.Range("B1").Value = 1.088173151016235
$ws.Range("C1").Value = 6.737033843994141
$ws.Range("D1").Value = 2.032927751541138
$ws.Range("E1").Value = 1.129818201065063
